# Refresh the cryptocurrency snapshot table on the active worksheet:
# updates "Price" (column D) and "Volume(1h)" (column E) for each coin row
# (rows 2-51) to the values from the latest scrape, per the scheduled
# GitHub Actions run. Numeric-looking price strings are written with a
# text number format first so Excel keeps them as literal text (matching
# the existing "22.450.56"-style values) instead of auto-converting them
# to numbers; the cell style is then reset to "Normal" so no stray
# per-cell number format lingers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.438.34"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.567.21"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3721"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3315"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.132"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07472"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.931"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.905"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").Value = "1.567.06"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001115"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06746"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.89%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.344"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("D24").Value = "22.428.10"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.389"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.561"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.020"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").Value = "1.739.61"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.012"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.116"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.730"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08319"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02454"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2270"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06392"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.369"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.285"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6285"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6137"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.773"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.209"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07219"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
